# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (positioned between "2020-Q4" and
#    "总计") carrying the per-fund holdings for the new quarter.
# 2. Update the "总计" (totals) worksheet with a new leading row for
#    2022-Q1, pushing the existing 2020-Q4 totals row down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by copying "总计" so it inherits the
# same header/row styling (bold header + border, style index used by the
# existing "总计" sheet), then insert it right before "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($total)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# Extend the styled header range (B1:D1) across to H1 so every header
# cell shares the same formatting as the existing "总计" header row.
$newSheet.Range("D1").Copy($newSheet.Range("E1:H1"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2: 长盛安睿一年持有混合A (012377)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'012377"
$newSheet.Range("C2").Value = "长盛安睿一年持有混合A"
$newSheet.Range("D2").Value = "'4.20"
$newSheet.Range("E2").Value = "'29.09"
$newSheet.Range("F2").Value = "'3.31"
$newSheet.Range("G2").Value = "'0.1390"
$newSheet.Range("H2").Value = 4

# Row 3: 长盛安睿一年持有混合C (012378)
$newSheet.Range("A2").Copy($newSheet.Range("A3"))
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'012378"
$newSheet.Range("C3").Value = "长盛安睿一年持有混合C"
$newSheet.Range("D3").Value = "'0.36"
$newSheet.Range("E3").Value = "'29.09"
$newSheet.Range("F3").Value = "'3.31"
$newSheet.Range("G3").Value = "'0.0119"
$newSheet.Range("H3").Value = 4

# The leading "'" above forces these numeric-looking strings to stay text
# (matching the source data's inline-string columns) instead of being
# coerced to numbers; drop the resulting quote-prefix text format so the
# cells fall back to the sheet's plain/default style, same as the source.
$newSheet.Range("B2").ClearFormats()
$newSheet.Range("D2:G2").ClearFormats()
$newSheet.Range("B3").ClearFormats()
$newSheet.Range("D3:G3").ClearFormats()

# ---------------------------------------------------------------------
# Step 2: update "总计" - push the existing 2020-Q4 row down to row 3 and
# write the new 2022-Q1 totals into row 2.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.15

# Restore the originally active sheet / selection.
$wb.Worksheets.Item("2020-Q4").Activate()
